# Update "Overview" sheet: roll the 5-period income statement window forward
# by one period (drop the oldest column's data, shift the remaining four
# periods left, and populate the newly-freed last column with the latest
# period's figures), per "update database and change read_price algorithm".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: period headers (column D..H) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (column D..H) ---
$ws.Range("D9").Value = "1399-05-12 (5)"
$ws.Range("E9").Value = "1400-04-15 (9)"
$ws.Range("F9").Value = "1401-04-05 (10)"
$ws.Range("G9").Value = "1402-02-27 (7)"

# H9's new text ("1402-02-27") looks like a date, so a plain .Value assignment
# would get auto-coerced into a date serial by Excel's input parser. Force
# text entry via a Text number format, then restore H9's original
# (unformatted) look by re-applying D9's cell format.
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-27"
$ws.Range("D9").Copy()
$ws.Range("H9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 11..27: shift D<-E<-F<-G<-H left, then fill H with the new period ---
# Each entry: row number, [newD, newE, newF, newG, newH]
$rows = @(
    @{ Row = 11; Vals = @(10793, 9310, 9939, 21388, 27661) },
    @{ Row = 12; Vals = @(-7479, -6985, -7781, -17518, -20351) },
    @{ Row = 13; Vals = @(3314, 2326, 2159, 3871, 7310) },
    @{ Row = 14; Vals = @(-453, -537, -411, -543, -1542) },
    @{ Row = 15; Vals = @("-", "-", "-", "-", "-") },
    @{ Row = 16; Vals = @("-", "-", 8, -561, -627) },
    @{ Row = 17; Vals = @(2861, 1789, 1755, 2767, 5141) },
    @{ Row = 18; Vals = @("-", "-", "-", -601, -1967) },
    @{ Row = 19; Vals = @(265, 331, 189, 627, 244) },
    @{ Row = 20; Vals = @(3126, 2120, 1944, 2792, 3418) },
    @{ Row = 21; Vals = @(-369, -131, -104, -216, -189) },
    @{ Row = 22; Vals = @(2757, 1989, 1840, 2576, 3229) },
    @{ Row = 23; Vals = @("-", "-", "-", "-", "-") },
    @{ Row = 24; Vals = @(2757, 1989, 1840, 2576, 3229) },
    @{ Row = 25; Vals = @(0, 0, 0, 0, 0) },
    @{ Row = 26; Vals = @(3056, 5457, 3096, 2653, 1984) },
    @{ Row = 27; Vals = @(0, 0, 0, 0, 0) }
)

$cols = @("D", "E", "F", "G", "H")

foreach ($entry in $rows) {
    $r = $entry.Row
    $vals = $entry.Vals
    for ($i = 0; $i -lt 5; $i++) {
        $addr = $cols[$i] + $r
        $ws.Range($addr).Value = $vals[$i]
    }
}
